$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = -7.456999999999999
$ws.Range("C10").Value = -13.1001
$ws.Range("C12").Value = -10.6672
$ws.Range("D15").Value = -7.968999999999999
$ws.Range("C18").Value = -12.91109999999999
$ws.Range("D20").Value = -7.5875
$ws.Range("D29").Value = -7.315099999999997
$ws.Range("D30").Value = -7.182400000000008
$ws.Range("D31").Value = -8.353499999999999
$ws.Range("C37").Value = -12.6477
$ws.Range("D40").Value = -7.910299999999995
$ws.Range("C55").Value = -13.74609999999999
$ws.Range("C68").Value = -11.56329999999999
$ws.Range("D68").Value = -7.130099999999995
$ws.Range("D76").Value = -7.372800000000001
$ws.Range("C77").Value = -12.26350000000001
$ws.Range("C78").Value = -12.32310000000001
$ws.Range("D87").Value = -7.866099999999997
$ws.Range("D88").Value = -7.405599999999996
$ws.Range("D96").Value = -7.615900000000005
$ws.Range("D98").Value = -8.408000000000005
$ws.Range("D101").Value = -7.750799999999999
$ws.Range("D102").Value = -7.800599999999995
